$wb = $excel.ActiveWorkbook

# 1. Rename the second sheet
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Include #0"

# 2. Update values on the Metadata sheet
$ws1 = $wb.Worksheets.Item(1)

# URL: pythia -> cicada
$ws1.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/ValueSet/eval-status"

# Date: updated timestamp
$ws1.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# Insert a new "Jurisdiction" row after "Contact" (row 10), pushing
# Description/Purpose/Copyright/Immutable down by one row.
$ws1.Rows.Item(11).Insert()

# Copy formatting from the row above so the new row matches the
# existing table styling (bordered/top-aligned/wrap style).
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)

$ws1.Range("A11").Value = "Jurisdiction"

# The Jurisdiction value is blank text (not merely an empty cell), so
# force text typing with a leading quote, then re-copy the row-above
# formatting to strip the resulting quote-prefix style flag while
# keeping the now-text empty value.
$ws1.Range("B11").Value = "'"
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)

# 3. Update the System URI on the second sheet: pythia -> cicada
$ws2.Range("B7").Value = "http://fhirfli.dev/fhir/ig/cicada/CodeSystem/EvalStatus"
